# Update "想去人数" (interest count) figures in the F column of the
# "展览" and "全部类型" worksheets, reflecting the newly generated
# output from the site rebuild (gh-pages @ 456a3b4).

$wb = $excel.ActiveWorkbook

# Map of row number -> new value for column F, identical on both sheets.
$updates = @{
    4  = 11715
    5  = 877
    6  = 117
    7  = 18
    10 = 173
    16 = 339
    17 = 1379
    18 = 81
    19 = 906
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
